# Add 5 additional (older) quarterly columns to the "Overview" income-statement sheet.
# Historically the sheet showed quarters D:H (1400/09 .. 1401/09). We now also need to
# show the 5 quarters preceding those, inserted at D:H, pushing the existing D:H data
# (and its formatting) to I:M.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 5 new columns at D, shifting the current D:H block (and its formatting) to I:M.
$ws.Range("D1:H1").EntireColumn.Insert()

# 2) Header row 8: period labels for the 5 newly-inserted (older) quarters.
$periods = @(
    "فصل دوم منتهی به 1399/06",
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06"
)
for ($i = 0; $i -lt $periods.Length; $i++) {
    $ws.Cells.Item(8, 4 + $i).Value = $periods[$i]
}

# 3) Header row 9: publish dates for the 5 newly-inserted quarters.
$pubDates = @(
    "1400-08-25 (4)",
    "1400-10-29 (2)",
    "1401-04-01 (8)",
    "1401-04-28 (2)",
    "1401-08-25 (4)"
)
for ($i = 0; $i -lt $pubDates.Length; $i++) {
    $ws.Cells.Item(9, 4 + $i).Value = $pubDates[$i]
}

# 4) Financial data rows 11-27, new D:H values (5 older quarters).
$data = @{
    11 = @(3780, 3305, 5026, 4742, 5390)
    12 = @(-2963, -2639, -3791, -3602, -4078)
    13 = @(817, 666, 1235, 1140, 1312)
    14 = @(-178, -171, -241, -260, -220)
    15 = @("-", "-", "-", "-", "-")
    16 = @(-33, -27, -18, -22, -18)
    17 = @(606, 468, 976, 858, 1074)
    18 = @(-1, -6, -10, -9, -7)
    19 = @(25, 1, 24, 2, 3)
    20 = @(630, 463, 990, 851, 1070)
    21 = @(-103, -88, -143, -160, -192)
    22 = @(527, 374, 847, 691, 878)
    23 = @("-", "-", "-", "-", "-")
    24 = @(527, 374, 847, 691, 878)
    25 = @(0, 0, 0, 0, 0)
    26 = @(20332, 16500, 18359, 19220, 17179)
    27 = @(0, 0, 0, 0, 0)
}
foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, 4 + $i).Value = $vals[$i]
    }
}

# 5) Restore the intended column widths for the newly-inserted D:H columns
#    (D=29, E=29, F=31, G=29, H=29 "characters", entered as ColumnWidth which
#    is offset by 5/6 from the stored OOXML column width).
$ws.Columns.Item(4).ColumnWidth = 29 - 5/6
$ws.Columns.Item(5).ColumnWidth = 29 - 5/6
$ws.Columns.Item(6).ColumnWidth = 31 - 5/6
$ws.Columns.Item(7).ColumnWidth = 29 - 5/6
$ws.Columns.Item(8).ColumnWidth = 29 - 5/6

# 6) Match the saved selection/cursor position.
$ws.Range("L6").Select()
